$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.978041291236877
$ws.Range("B1").Value = 2.152907848358154
$ws.Range("C1").Value = 2.44899320602417
$ws.Range("D1").Value = 3.743246078491211
$ws.Range("E1").Value = 1.31428861618042
